$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "synthetic"
$ws.Range("C8").Value = "original"
$ws.Range("D8").Value = $true
$ws.Range("E8").Value = "C:\Users\franz\Documents\work\projects\arp\data\synthetic_data\synthetic_data_original_textured_unclipped_vtp_paraview"
$ws.Range("F8").Value = ".vtp"
$ws.Range("G8").Value = 200
$ws.Range("H8").Value = "M"
$ws.Range("I8").Value = $false
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = "automatic"
$ws.Range("M8").Value = $false
$ws.Range("N8").Value = $true
$ws.Range("O8").Value = $true

$ws.Range("L11").Select() | Out-Null
